# Adds a check to see if positive cost resources are needed for the RPS.
# CCS plants (hard coal w CCS, natural gas combined cycle w CCS, biomass w CCS,
# lignite w CCS -- rows 19-22) are moved to guaranteed dispatch, i.e. their
# BAU Guaranteed Dispatch Percentage is switched from 0 to 1 across all years.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# Bring the BGDPbES sheet to the front / make it the active tab.
$ws.Activate()

# Flip the guaranteed-dispatch flag to 1 for the CCS rows (19-22), all years (B:AK).
$ws.Range("B19:AK22").Value = 1

# Match the saved selection/active cell on the sheet.
$ws.Range("B19:AK22").Select() | Out-Null
